$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the projection numbers for the players that stay on the board.
$ws.Range("B2").Value = 52
$ws.Range("C2").Value = 0

$ws.Range("B3").Value = 52
$ws.Range("C3").Value = 48

$ws.Range("B4").Value = 43
$ws.Range("C4").Value = 0

$ws.Range("B5").Value = 43
$ws.Range("C5").Value = 0

$ws.Range("B6").Value = 39
$ws.Range("C6").Value = 40

$ws.Range("B7").Value = 37
$ws.Range("C7").Value = 34

$ws.Range("B8").Value = 33
$ws.Range("C8").Value = 29

$ws.Range("B9").Value = 27
$ws.Range("C9").Value = 28

$ws.Range("B10").Value = 23
$ws.Range("C10").Value = 14

$ws.Range("B11").Value = 24
$ws.Range("C11").Value = 29

$ws.Range("B12").Value = 26
$ws.Range("C12").Value = 23

$ws.Range("B13").Value = 14
$ws.Range("C13").Value = 15

$ws.Range("B14").Value = 15
$ws.Range("C14").Value = 1

$ws.Range("B15").Value = 16
$ws.Range("C15").Value = 22

$ws.Range("B16").Value = 14
$ws.Range("C16").Value = 2

$ws.Range("B17").Value = 12
$ws.Range("C17").Value = 16

# Drop MOUSSA DIABATE from the board: row 18 becomes the next player
# (DEVONTE GRAHAM), and everyone below moves up one spot.
$ws.Range("A18").Value = "DEVONTE GRAHAM"
$ws.Range("B18").Value = 14
$ws.Range("C18").Value = 10

$ws.Range("A19").Value = "ALEX LEN"
$ws.Range("B19").Value = 6
$ws.Range("C19").Value = 1
$ws.Range("A19").Font.Bold = $false

$ws.Range("A20").Value = "JEVON CARTER"
$ws.Range("B20").Value = 7
$ws.Range("C20").Value = 16

$ws.Range("A21").Value = "MATISSE THYBULLE"
$ws.Range("B21").Value = 7
$ws.Range("C21").Value = 11

$ws.Range("A22").Value = "MARK WILLIAMS"
$ws.Range("B22").Value = 6
$ws.Range("C22").Value = 0

# Add the new pick, FRED VANVLEET, as the last row of the board.
$ws.Range("A23").Value = "FRED VANVLEET"
$ws.Range("B23").Value = 38
$ws.Range("C23").Value = 39

# Move the active selection to B12, as on the author's screen when they
# finished making these picks.
$ws.Range("B12").Select()
